# fix: revert admin dev default; seed customers only when table empty;
# autosave on customer select when hours/day present
#
# Helper: VBA-style RGB() -> packed BGR-ordered OLE color int Excel expects.
function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Timesheet")

# ---------------------------------------------------------------------
# 1. Weekly Timesheet: worked rows reset to the 8h/day admin default,
#    with Rate/Total zeroed out (pending re-entry), and new employee
#    names for this period.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Evans"
$ws.Range("B3").Value = "Oglesby"
$ws.Range("B4").Value = "Muncey"
$ws.Range("B5").Value = "Lucas"
$ws.Range("B6").Value = "Bailey"

$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 8
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 8

$ws.Range("E2:E6").Value = 0
$ws.Range("F2:F6").Value = 0

# SUBTOTAL row: 40 regular hours, total reverts to 0 until re-rated.
$ws.Range("C8").Value = 40
$ws.Range("F8").Value = 0
$ws.Range("D8").Value = "Reg: 40 / OT: 0"

# ---------------------------------------------------------------------
# 2. New summary rows 11-13 under the existing table.
# ---------------------------------------------------------------------
$hourly = $ws.Range("A11:F11")
$hourly.Font.Bold = $true
$hourly.Interior.Color = (RGB 250 243 224)
$ws.Range("A11").Value = "HOURLY SUBTOTAL"
$ws.Range("F11").NumberFormat = '"$"#,##0.00'
$ws.Range("F11").Value = 0

$admin = $ws.Range("A12:F12")
$admin.Font.Bold = $true
$admin.Interior.Color = (RGB 250 243 224)
$ws.Range("A12").Value = "ADMIN SUBTOTAL"
$ws.Range("F12").NumberFormat = '"$"#,##0.00'
$ws.Range("F12").Value = 0

$grand = $ws.Range("A13:F13")
$grand.Font.Bold = $true
$grand.Font.Color = (RGB 255 0 0)
$grand.Interior.Color = (RGB 232 248 224)
$ws.Range("A13").Value = "GRAND TOTAL"
$ws.Range("F13").NumberFormat = '"$"#,##0.00'
$ws.Range("F13").Value = 0

# ---------------------------------------------------------------------
# 3. New "Jason Schema" sheet - flattened per-employee/per-day export,
#    seeded from the (now-empty) timesheet rows above.
# ---------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item(1)
$js = $wb.Worksheets.Add($null, $sheet1)
$js.Name = "Jason Schema"

$js.Columns.Item(1).ColumnWidth = 20
$js.Columns.Item(2).ColumnWidth = 18
$js.Columns.Item(3).ColumnWidth = 12
$js.Columns.Item(4).ColumnWidth = 25
$js.Columns.Item(5).ColumnWidth = 8
$js.Columns.Item(6).ColumnWidth = 10
$js.Columns.Item(7).ColumnWidth = 12
$js.Columns.Item(8).ColumnWidth = 10
$js.Columns.Item(9).ColumnWidth = 30

$jsHeader = $js.Range("A1:I1")
$jsHeader.Font.Bold = $true

$js.Range("A1").Value = "Employee"
$js.Range("B1").Value = "Employee ID"
$js.Range("C1").Value = "Date"
$js.Range("D1").Value = "Client"
$js.Range("E1").Value = "Hours"
$js.Range("F1").Value = "Rate"
$js.Range("G1").Value = "Total"
$js.Range("H1").Value = "Type"
$js.Range("I1").Value = "Notes"

$js.Range("F1:G1").NumberFormat = '"$"#,##0.00'
$js.Range("F2:G6").NumberFormat = '"$"#,##0.00'

$employee = "Boban Abbate"
$employeeId = "emp_pw6be4hd"

# Leading "'" forces text storage so the YYYY-MM-DD strings round-trip as
# shared-string text instead of being coerced into date serials.
$dates = @("'2026-01-05", "'2026-01-06", "'2026-01-07", "'2026-01-08", "'2026-01-09")
$clients = @("Evans", "Oglesby", "Muncey", "Lucas", "Bailey")

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $js.Cells.Item($r, 1).Value = $employee
    $js.Cells.Item($r, 2).Value = $employeeId
    $js.Cells.Item($r, 3).Value = $dates[$i]
    $js.Cells.Item($r, 4).Value = $clients[$i]
    $js.Cells.Item($r, 5).Value = 8
    $js.Cells.Item($r, 6).Value = 0
    $js.Cells.Item($r, 7).Value = 0
    $js.Cells.Item($r, 8).Value = "Regular"
    $js.Cells.Item($r, 9).Value = ""
}

$ws.Activate()
